$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove whole paragraphs that disappear in the target revision.
#    Deleted bottom-up so earlier paragraph indices stay stable.
# ------------------------------------------------------------------

# "डीकन, डोसेटिज़्म" (italic) paragraph right after the "ड" Heading2
$d.Paragraphs.Item(11).Range.Delete()

# "This PDF version is provided under the same license." paragraph
$d.Paragraphs.Item(6).Range.Delete()

# "License Information" Heading2 paragraph
$d.Paragraphs.Item(4).Range.Delete()

# ------------------------------------------------------------------
# 2. Rewrite the big license/attribution paragraph (now paragraph 4).
# ------------------------------------------------------------------

# Replace the big run-set "Biblica Bible Dictionary" ... "CC BY-SA 4.0 license"
# (which spans two hyperlinks) through the trailing "." with a single new run.
$rngStart = $d.Content
$rngStart.Find.Execute("Biblica Bible Dictionary", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos = $rngStart.Start

$rngEnd = $d.Content
$rngEnd.Find.Execute("CC BY-SA 4.0 license", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $rngEnd.End + 1  # include the trailing "."

$bigRange = $d.Range($startPos, $endPos)
$bigRange.Text = " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."

# Replace the bold heading-style run inside this paragraph.
$d.Content.Find.Execute("मुख्य शब्द (Biblica)", $true, $false, $false, $false, $false, $true, 1, $false, "Biblica Study Notes (Key Terms)", 2)

# Replace " (Hindi) is based on" -> " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. "
$d.Content.Find.Execute(" (Hindi) is based on", $true, $false, $false, $false, $false, $true, 1, $false, " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. ", 2)

# Replace ": " (the run immediately before the big merged run) with "Biblica Study Notes"
$rngColon = $d.Content
$rngColon.Find.Execute(" has been adapted", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$colonTarget = $d.Range($rngColon.Start - 2, $rngColon.Start)
$colonTarget.Text = "Biblica Study Notes"
